$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")
$wsPlaylist = $wb.Worksheets.Item("Add_Playlist")

# --- Add new row 3 to the "Search" sheet ("add play list" test case) ---

# A3 is a number-looking label ("2") that must be stored as text (like A2 = "1"),
# so force text format before assigning the value, then copy the cell format
# from A2 (which already carries the correct quote-prefixed text style).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("B3").Value = "Nhập đầy đủ tên một bài hát "
$ws.Range("C3").Value = "Nơi này có anh"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Hiển thị bài hát nơi này có anh trên trang web"
$ws.Range("F3").Value = "Pass"
$ws.Range("G3").Value = "//div[@class='sm7ZnbOO1Zfg9cupYgPN']/a[@title='Nơi Này Có Anh']"

$ws.Rows("3").RowHeight = 45

# Update the sheet's selection to reflect the newly added row.
$ws.Range("F2:F3").Select()

# --- Update the "Add_Playlist" sheet's stored selection ---
$wsPlaylist.Range("A1:G1").Select()

Write-Output "done"
